# Apply updates to the "Notes" sheet of the Uganda WASH performance score workbook.
#
# Summary of changes (per commit "user data and updated uganda financial calcs"):
#  - Description text reworded.
#  - Source text reworded, and a new "Source-link" line added right after it.
#  - License text reworded, and a new "More information on licensing..." line
#    added right after it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notes")

# --- Source block (rows 1-4 stay put; insert a new row 5 for the source link) ---
$ws.Rows.Item(5).Insert()
$ws.Range("A5").Value = "Source-link: http://www.mwe.go.ug/index.php?option=com_docman&task=cat_view&Itemid=223&gid=15"

$ws.Range("A2").Value = "Description: District Wash Performance Score"
$ws.Range("A4").Value = "Source: Water and Environment Sector Performance Reports 2010-2014 - Ministry of Water and Environment"

# --- Licensing block: after the insert above, the old license row (13) is now
#     row 14. Insert a new row 15 for the extra licensing-info line. ---
$ws.Rows.Item(15).Insert()
$ws.Range("A14").Value = "It is licensed under a Creative Commons Attribution 4.0 International license."
$ws.Range("A15").Value = "More information on licensing is available here: https://creativecommons.org/licenses/by/4.0/"
